# Ajout de plusieurs années
# - Adds a new worksheet "INFOS-PARCOURS" at the end of the workbook
# - Populates it with a small header/value table
# - Makes the new sheet the active one (tabSelected moves to it)

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "INFOS-PARCOURS"

# Header row
$ws3.Range("A1").Value = "id_parcour"
$ws3.Range("B1").Value = "id_semestre"

# Data row
$ws3.Range("A2").Value = "RATIO"
$ws3.Range("B2").Value = "6_2023-2024"

# Column B is a bit wider than default
$ws3.Columns.Item(2).ColumnWidth = 12.1666666666667

# Make the new sheet the active tab/selection, matching the saved view state
$ws3.Activate()
$ws3.Range("I9").Select()
